# Auto-generated edit script applying the Moogle_Profits market-data refresh.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# row 4 (@@ -828,22 +828,22 @@)
$ws.Range("H4").Value = 962.5
$ws.Range("I4").Value = 453.5
$ws.Range("K4").Value = 453.5
$ws.Range("M4").Value = -339.5

# row 28 (@@ -2004,25 +2004,25 @@)
$ws.Range("H28").Value = 593.6875
$ws.Range("I28").Value = 277.3
$ws.Range("J28").Value = 1121
$ws.Range("K28").Value = 277.3
$ws.Range("L28").Value = 1121
$ws.Range("M28").Value = 207.7
$ws.Range("N28").Value = -2091

# row 40 (@@ -2607,25 +2607,25 @@)
$ws.Range("H40").Value = 3845.8572
$ws.Range("J40").Value = 4031.5454
$ws.Range("L40").Value = 4031.5454
$ws.Range("N40").Value = -4381.5454

# row 62 (@@ -3712,25 +3712,25 @@)
$ws.Range("H62").Value = 8972
$ws.Range("J62").Value = 14444
$ws.Range("L62").Value = 14444
$ws.Range("N62").Value = -15692

# row 65 (@@ -3865,25 +3865,25 @@)
$ws.Range("H65").Value = 8972
$ws.Range("J65").Value = 14444
$ws.Range("L65").Value = 72220
$ws.Range("N65").Value = -78460

# row 96 (@@ -5447,25 +5447,25 @@)
$ws.Range("H96").Value = 907.6429000000001
$ws.Range("I96").Value = 688.5714
$ws.Range("J96").Value = 1126.7142
$ws.Range("K96").Value = 2065.7142
$ws.Range("L96").Value = 3380.1426
$ws.Range("M96").Value = -692.7142000000003
$ws.Range("N96").Value = -6126.142599999999

# row 100 (@@ -5652,25 +5652,25 @@)
$ws.Range("H100").Value = 2096.375
$ws.Range("I100").Value = 2253.1667
$ws.Range("J100").Value = 2002.3
$ws.Range("K100").Value = 2253.1667
$ws.Range("L100").Value = 2002.3
$ws.Range("M100").Value = -1712.1667
$ws.Range("N100").Value = -3084.3

# row 129 (@@ -7103,25 +7103,25 @@)
$ws.Range("H129").Value = 3157.625
$ws.Range("I129").Value = 3077.5217
$ws.Range("J129").Value = 5000
$ws.Range("K129").Value = 9232.5651
$ws.Range("L129").Value = 15000
$ws.Range("M129").Value = -4232.5651
$ws.Range("N129").Value = -25000

# row 137 (@@ -7495,25 +7495,25 @@)
$ws.Range("H137").Value = 3805.2
$ws.Range("I137").Value = 3196.366
$ws.Range("J137").Value = 4539.3823
$ws.Range("K137").Value = 9589.098
$ws.Range("L137").Value = 13618.1469
$ws.Range("M137").Value = -7039.098
$ws.Range("N137").Value = -18718.1469

# row 138 (@@ -7547,25 +7547,25 @@)
$ws.Range("H138").Value = 13307779
$ws.Range("I138").Value = 2068669.6
$ws.Range("J138").Value = 26321486
$ws.Range("K138").Value = 6206008.800000001
$ws.Range("L138").Value = 78964458
$ws.Range("M138").Value = -6200868.800000001
$ws.Range("N138").Value = -78974738

$ws = $wb.Worksheets.Item("ARM")
# row 32 (@@ -9310,22 +9310,25 @@)
$ws.Range("H32").Value = 3809.1052
$ws.Range("I32").Value = 2243.236
$ws.Range("J32").Value = 31994.75
$ws.Range("K32").Value = 2243.236
$ws.Range("L32").Value = 31994.75
$ws.Range("N32").Value = -32568.75
$ws.Range("M32").Value = -1956.236

# row 61 (@@ -10713,25 +10716,25 @@)
$ws.Range("H61").Value = 8190.2856
$ws.Range("I61").Value = 3291.6667
$ws.Range("J61").Value = 14721.777
$ws.Range("K61").Value = 3291.6667
$ws.Range("L61").Value = 14721.777
$ws.Range("M61").Value = -3079.6667
$ws.Range("N61").Value = -15145.777

# row 74 (@@ -11347,25 +11350,25 @@)
$ws.Range("H74").Value = 2942.1462
$ws.Range("I74").Value = 1733.52
$ws.Range("J74").Value = 4830.625
$ws.Range("K74").Value = 1733.52
$ws.Range("L74").Value = 4830.625
$ws.Range("M74").Value = -859.52
$ws.Range("N74").Value = -6578.625

# row 77 (@@ -11497,25 +11500,25 @@)
$ws.Range("H77").Value = 2942.1462
$ws.Range("I77").Value = 1733.52
$ws.Range("J77").Value = 4830.625
$ws.Range("K77").Value = 8667.6
$ws.Range("L77").Value = 24153.125
$ws.Range("M77").Value = -4299.6
$ws.Range("N77").Value = -32889.125

# row 110 (@@ -13138,25 +13141,25 @@)
$ws.Range("H110").Value = 3304
$ws.Range("I110").Value = 3999
$ws.Range("J110").Value = 2956.5
$ws.Range("K110").Value = 3999
$ws.Range("L110").Value = 2956.5
$ws.Range("M110").Value = -1954
$ws.Range("N110").Value = -7046.5

# row 122 (@@ -13723,22 +13726,22 @@)
$ws.Range("H122").Value = 2443.875
$ws.Range("I122").Value = 2592.9048
$ws.Range("K122").Value = 7778.714399999999
$ws.Range("M122").Value = -5328.714399999999

# row 136 (@@ -14394,25 +14397,25 @@)
$ws.Range("H136").Value = 8190.2856
$ws.Range("I136").Value = 3291.6667
$ws.Range("J136").Value = 14721.777
$ws.Range("K136").Value = 9875.000100000001
$ws.Range("L136").Value = 44165.331
$ws.Range("M136").Value = -7325.000100000001
$ws.Range("N136").Value = -49265.331

$ws = $wb.Worksheets.Item("BSM")
# row 99 (@@ -19535,25 +19538,25 @@)
$ws.Range("H99").Value = 2309.9
$ws.Range("I99").Value = 1571.4286
$ws.Range("J99").Value = 4033
$ws.Range("K99").Value = 1571.4286
$ws.Range("L99").Value = 4033
$ws.Range("M99").Value = -73.42859999999996
$ws.Range("N99").Value = -7029

# row 105 (@@ -19835,25 +19838,25 @@)
$ws.Range("H105").Value = 9219.559999999999
$ws.Range("I105").Value = 7499.591
$ws.Range("J105").Value = 21832.666
$ws.Range("K105").Value = 7499.591
$ws.Range("L105").Value = 21832.666
$ws.Range("M105").Value = -5752.591
$ws.Range("N105").Value = -25326.666

$ws = $wb.Worksheets.Item("CRP")
# row 22 (@@ -22710,25 +22713,25 @@)
$ws.Range("H22").Value = 1810.8125
$ws.Range("J22").Value = 1823.5834
$ws.Range("L22").Value = 1823.5834
$ws.Range("N22").Value = -2523.5834

# row 58 (@@ -24477,25 +24480,25 @@)
$ws.Range("H58").Value = 7890.0527
$ws.Range("I58").Value = 4962.4614
$ws.Range("J58").Value = 14233.167
$ws.Range("K58").Value = 4962.4614
$ws.Range("L58").Value = 14233.167
$ws.Range("M58").Value = -4759.4614
$ws.Range("N58").Value = -14639.167

# row 86 (@@ -25876,25 +25879,25 @@)
$ws.Range("H86").Value = 7686.375
$ws.Range("I86").Value = 7098.6
$ws.Range("J86").Value = 8666
$ws.Range("K86").Value = 7098.6
$ws.Range("L86").Value = 8666
$ws.Range("M86").Value = -5975.6
$ws.Range("N86").Value = -10912

# row 89 (@@ -26026,25 +26029,25 @@)
$ws.Range("H89").Value = 7686.375
$ws.Range("I89").Value = 7098.6
$ws.Range("J89").Value = 8666
$ws.Range("K89").Value = 35493
$ws.Range("L89").Value = 43330
$ws.Range("M89").Value = -29877
$ws.Range("N89").Value = -54562

# row 99 (@@ -26522,25 +26525,25 @@)
$ws.Range("H99").Value = 2384.175
$ws.Range("I99").Value = 1577.2812
$ws.Range("J99").Value = 5611.75
$ws.Range("K99").Value = 1577.2812
$ws.Range("L99").Value = 5611.75
$ws.Range("M99").Value = -79.2811999999999
$ws.Range("N99").Value = -8607.75

# row 105 (@@ -26822,22 +26825,22 @@)
$ws.Range("H105").Value = 1867.3334
$ws.Range("I105").Value = 1972.5714
$ws.Range("K105").Value = 1972.5714
$ws.Range("M105").Value = -225.5714

# row 126 (@@ -27863,25 +27866,25 @@)
$ws.Range("H126").Value = 2384.175
$ws.Range("I126").Value = 1577.2812
$ws.Range("J126").Value = 5611.75
$ws.Range("K126").Value = 4731.8436
$ws.Range("L126").Value = 16835.25
$ws.Range("M126").Value = -2261.8436
$ws.Range("N126").Value = -21775.25

# row 132 (@@ -28154,22 +28157,22 @@)
$ws.Range("H132").Value = 4810.88
$ws.Range("I132").Value = 4070.7222
$ws.Range("K132").Value = 12212.1666
$ws.Range("M132").Value = -9682.1666

# row 134 (@@ -28255,22 +28258,22 @@)
$ws.Range("H134").Value = 9301.75
$ws.Range("I134").Value = 5636.6665
$ws.Range("K134").Value = 16909.9995
$ws.Range("M134").Value = -14374.9995

# row 136 (@@ -28356,25 +28359,25 @@)
$ws.Range("H136").Value = 7890.0527
$ws.Range("I136").Value = 4962.4614
$ws.Range("J136").Value = 14233.167
$ws.Range("K136").Value = 14887.3842
$ws.Range("L136").Value = 42699.501
$ws.Range("M136").Value = -12337.3842
$ws.Range("N136").Value = -47799.501

$ws = $wb.Worksheets.Item("CUL")
# row 80 (@@ -32686,25 +32689,25 @@)
$ws.Range("H80").Value = 4135.467
$ws.Range("I80").Value = 4500
$ws.Range("J80").Value = 4079.3845
$ws.Range("K80").Value = 13500
$ws.Range("L80").Value = 12238.1535
$ws.Range("M80").Value = -12564
$ws.Range("N80").Value = -14110.1535

# row 83 (@@ -32839,25 +32842,25 @@)
$ws.Range("H83").Value = 4135.467
$ws.Range("I83").Value = 4500
$ws.Range("J83").Value = 4079.3845
$ws.Range("K83").Value = 40500
$ws.Range("L83").Value = 36714.4605
$ws.Range("M83").Value = -35820
$ws.Range("N83").Value = -46074.4605

# row 113 (@@ -34348,25 +34351,25 @@)
$ws.Range("H113").Value = 954.125
$ws.Range("J113").Value = 1068.75
$ws.Range("L113").Value = 3206.25
$ws.Range("N113").Value = -7546.25

$ws = $wb.Worksheets.Item("GSM")
# row 2 (@@ -35920,25 +35923,25 @@)
$ws.Range("H2").Value = 1420.9524
$ws.Range("I2").Value = 927.94116
$ws.Range("J2").Value = 3516.25
$ws.Range("K2").Value = 927.94116
$ws.Range("L2").Value = 3516.25
$ws.Range("M2").Value = -814.94116
$ws.Range("N2").Value = -3742.25

# row 10 (@@ -36318,25 +36321,25 @@)
$ws.Range("H10").Value = 55160.8
$ws.Range("J10").Value = 85300.336
$ws.Range("L10").Value = 85300.336
$ws.Range("N10").Value = -85638.336

# row 97 (@@ -40605,22 +40608,22 @@)
$ws.Range("H97").Value = 955.4783
$ws.Range("I97").Value = 920.5625
$ws.Range("K97").Value = 920.5625
$ws.Range("M97").Value = -424.5625

# row 102 (@@ -40853,25 +40856,25 @@)
$ws.Range("H102").Value = 3699.923
$ws.Range("I102").Value = 1833.3334
$ws.Range("J102").Value = 7899.75
$ws.Range("K102").Value = 1833.3334
$ws.Range("L102").Value = 7899.75
$ws.Range("M102").Value = -211.3334
$ws.Range("N102").Value = -11143.75

# row 113 (@@ -41398,22 +41401,22 @@)
$ws.Range("H113").Value = 6040.091
$ws.Range("I113").Value = 4583.125
$ws.Range("K113").Value = 4583.125
$ws.Range("M113").Value = -2413.125

$ws = $wb.Worksheets.Item("LTW")
# row 7 (@@ -43167,22 +43170,22 @@)
$ws.Range("H7").Value = 109789.7
$ws.Range("I7").Value = 210979.6
$ws.Range("K7").Value = 210979.6
$ws.Range("M7").Value = -210867.6

# row 11 (@@ -43369,22 +43372,22 @@)
$ws.Range("H11").Value = 17499.5
$ws.Range("I11").Value = 0
$ws.Range("J11").Value = 17499.5
$ws.Range("K11").Value = 0
$ws.Range("L11").Value = 17499.5
$ws.Range("M11").ClearContents()
$ws.Range("N11").Value = -17779.5

# row 82 (@@ -46863,22 +46866,22 @@)
$ws.Range("H82").Value = 2769.6316
$ws.Range("I82").Value = 5833.2
$ws.Range("K82").Value = 5833.2
$ws.Range("M82").Value = -5472.2

# row 85 (@@ -47013,22 +47016,22 @@)
$ws.Range("H85").Value = 2769.6316
$ws.Range("I85").Value = 5833.2
$ws.Range("K85").Value = 5833.2
$ws.Range("M85").Value = -4585.2

# row 100 (@@ -47739,22 +47742,22 @@)
$ws.Range("H100").Value = 5750.1113
$ws.Range("I100").Value = 6250.125
$ws.Range("K100").Value = 6250.125
$ws.Range("M100").Value = -5709.125

# row 126 (@@ -49004,22 +49007,22 @@)
$ws.Range("H126").Value = 109789.7
$ws.Range("I126").Value = 210979.6
$ws.Range("K126").Value = 632938.8
$ws.Range("M126").Value = -630468.8

# row 132 (@@ -49292,25 +49295,25 @@)
$ws.Range("H132").Value = 6334.222
$ws.Range("I132").Value = 2694.8462
$ws.Range("J132").Value = 15796.6
$ws.Range("K132").Value = 8084.5386
$ws.Range("L132").Value = 47389.8
$ws.Range("M132").Value = -5554.5386
$ws.Range("N132").Value = -52449.8

$ws = $wb.Worksheets.Item("WVR")
# row 100 (@@ -54675,25 +54678,25 @@)
$ws.Range("H100").Value = 866.4286
$ws.Range("I100").Value = 647.5
$ws.Range("J100").Value = 1158.3334
$ws.Range("K100").Value = 1295
$ws.Range("L100").Value = 2316.6668
$ws.Range("M100").Value = -754
$ws.Range("N100").Value = -3398.6668

# row 107 (@@ -55021,25 +55024,25 @@)
$ws.Range("H107").Value = 2468.75
$ws.Range("I107").Value = 2035.2222
$ws.Range("J107").Value = 3026.1428
$ws.Range("K107").Value = 6105.6666
$ws.Range("L107").Value = 9078.428400000001
$ws.Range("M107").Value = -4185.6666
$ws.Range("N107").Value = -12918.4284

# row 126 (@@ -55961,22 +55964,22 @@)
$ws.Range("H126").Value = 2847.625
$ws.Range("I126").Value = 3155.2
$ws.Range("K126").Value = 9465.599999999999
$ws.Range("M126").Value = -6995.599999999999
